$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-change (E) columns per the latest crypto data snapshot.

$ws.Range('D2').Value = '27.257.24'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '1.786.54'
$ws.Range('E3').Value = '  -1.59%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = "'333.62"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.24%  '
$ws.Range('E7').Value = '  -1.52%  '
$ws.Range('D8').Value = "'48.72"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.22%  '
$ws.Range('E9').Value = '  -2.97%  '
$ws.Range('E10').Value = '  -3.39%  '
$ws.Range('D11').Value = "'0.07488"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').Value = "'21.91"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('E14').Value = '  -2.49%  '
$ws.Range('D15').Value = '1.789.63'
$ws.Range('E15').Value = '  -1.37%  '
$ws.Range('D16').Value = "'7.087"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.15%  '
$ws.Range('D17').Value = "'0.00001096"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.49%  '
$ws.Range('E18').Value = '  -2.17%  '
$ws.Range('D19').Value = "'83.88"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.77%  '
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('D21').Value = "'6.638"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.32%  '
$ws.Range('E22').Value = '  -2.84%  '
$ws.Range('D23').Value = '27.268.17'
$ws.Range('E23').Value = '  -0.60%  '
$ws.Range('E24').Value = '  -6.48%  '
$ws.Range('D25').Value = "'2.417"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.44%  '
$ws.Range('D26').Value = "'1.508"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.86%  '
$ws.Range('D27').Value = "'2.542"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.53%  '
$ws.Range('D28').Value = "'21.32"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.22%  '
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('D30').Value = '1.990.70'
$ws.Range('E30').Value = '  -1.34%  '
$ws.Range('D31').Value = "'134.16"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.03%  '
$ws.Range('D32').Value = "'4.015"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.25%  '
$ws.Range('D33').Value = "'6.088"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.81%  '
$ws.Range('D34').Value = "'0.08696"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.58%  '
$ws.Range('D35').Value = "'13.30"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.88%  '
$ws.Range('D36').Value = "'1.658"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.81%  '
$ws.Range('D37').Value = "'0.6954"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.66%  '
$ws.Range('D38').Value = "'5.452"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.52%  '
$ws.Range('E39').Value = '  -3.00%  '
$ws.Range('D40').Value = "'8.820"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('D41').Value = "'0.06322"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.75%  '
$ws.Range('D42').Value = "'0.02338"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.32%  '
$ws.Range('E43').Value = '  -1.81%  '
$ws.Range('D44').Value = "'14.45"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.46%  '
$ws.Range('D45').Value = "'0.6514"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.42%  '
$ws.Range('D47').Value = "'3.835"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.66%  '
$ws.Range('D48').Value = "'2.153"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.72%  '
$ws.Range('E49').Value = '  -2.70%  '
$ws.Range('D50').Value = "'0.07129"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.23%  '
$ws.Range('D51').Value = "'79.13"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.07%  '
